$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = "60.445.30"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "2.372.06"
$ws.Range("E3").Value = "  -3.52%  "

# Row 4: update D4, E4
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: update D5, E5
$ws.Range("D5").Formula = "'540.50"
$ws.Range("E5").Value = "  -0.78%  "

# Row 6: update D6, E6
$ws.Range("D6").Formula = "'139.36"
$ws.Range("E6").Value = "  -3.28%  "

# Row 7: update D7, E7
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8: update D8, E8
$ws.Range("D8").Formula = "'0.575"
$ws.Range("E8").Value = "  -5.99%  "

# Row 9: update D9, E9
$ws.Range("D9").Value = "2.366.33"
$ws.Range("E9").Value = "  -3.48%  "

# Row 10: update E10
$ws.Range("E10").Value = "  -1.23%  "

# Row 11: update E11
$ws.Range("E11").Value = "  +0.48%  "

# Row 12: update E12
$ws.Range("E12").Value = "  -0.40%  "

# Row 13: update D13, E13
$ws.Range("D13").Formula = "'0.341"
$ws.Range("E13").Value = "  -2.72%  "

# Row 14: update E14
$ws.Range("E14").Value = "  -2.18%  "

# Row 15: update D15, E15
$ws.Range("D15").Value = "2.794.09"
$ws.Range("E15").Value = "  -3.72%  "

# Row 16: update D16, E16
$ws.Range("D16").Formula = "'0.0000163"
$ws.Range("E16").Value = "  +0.78%  "

# Row 17: update D17, E17
$ws.Range("D17").Value = "60.002.89"
$ws.Range("E17").Value = "  -1.21%  "

# Row 18: update D18, E18
$ws.Range("D18").Value = "2.369.92"
$ws.Range("E18").Value = "  -4.06%  "

# Row 19: update D19, E19
$ws.Range("D19").Formula = "'10.55"
$ws.Range("E19").Value = "  -4.27%  "

# Row 20: update D20, E20
$ws.Range("D20").Formula = "'316.51"
$ws.Range("E20").Value = "  -0.27%  "

# Row 21: update D21, E21
$ws.Range("D21").Formula = "'4.08"
$ws.Range("E21").Value = "  -1.56%  "

# Row 22: update E22
$ws.Range("E22").Value = "  -3.60%  "

# Row 23: update E23
$ws.Range("E23").Value = "  -0.20%  "

# Row 24: update D24, E24
$ws.Range("D24").Formula = "'1.80"
$ws.Range("E24").Value = "  +4.70%  "

# Row 25: update E25
$ws.Range("E25").Value = "  -0.22%  "

# Row 26: update E26
$ws.Range("E26").Value = "  -0.24%  "

# Row 27: update D27, E27
$ws.Range("D27").Value = "2.482.97"
$ws.Range("E27").Value = "  -4.84%  "

# Row 28: update D28, E28
$ws.Range("D28").Value = "0.0₃0922"
$ws.Range("E28").Value = "  -3.67%  "

# Row 29: update D29, E29
$ws.Range("D29").Formula = "'7.67"
$ws.Range("E29").Value = "  +1.36%  "

# Row 30: update D30, E30
$ws.Range("D30").Formula = "'519.16"
$ws.Range("E30").Value = "  -1.55%  "

# Row 31: update B31, C31, D31, E31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Formula = "'1.41"
$ws.Range("E31").Value = "  -4.15%  "

# Row 32: update B32, C32, D32, E32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Formula = "'7.95"
$ws.Range("E32").Value = "  -3.52%  "

# Row 33: update E33
$ws.Range("E33").Value = "  -2.52%  "

# Row 34: update E34
$ws.Range("E34").Value = "  -3.39%  "

# Row 35: update E35
$ws.Range("E35").Value = "  +0.33%  "

# Row 36: update D36, E36
$ws.Range("D36").Formula = "'0.998"
$ws.Range("E36").Value = "  -0.07%  "

# Row 37: update D37, E37
$ws.Range("D37").Formula = "'5.42"
$ws.Range("E37").Value = "  -5.58%  "

# Row 38: update D38, E38
$ws.Range("D38").Formula = "'4.61"
$ws.Range("E38").Value = "  -4.04%  "

# Row 39: update D39, E39
$ws.Range("D39").Formula = "'0.373"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40: update D40, E40
$ws.Range("D40").Formula = "'17.98"
$ws.Range("E40").Value = "  -1.50%  "

# Row 41: update E41
$ws.Range("E41").Value = "  -0.08%  "

# Row 42: update E42
$ws.Range("E42").Value = "  +1.81%  "

# Row 43: update D43, E43
$ws.Range("D43").Formula = "'137.02"
$ws.Range("E43").Value = "  -4.23%  "

# Row 44: update D44, E44
$ws.Range("D44").Formula = "'40.09"
$ws.Range("E44").Value = "  +0.22%  "

# Row 45: update D45, E45
$ws.Range("D45").Formula = "'2.17"
$ws.Range("E45").Value = "  -4.17%  "

# Row 46: update D46, E46
$ws.Range("D46").Formula = "'139.19"
$ws.Range("E46").Value = "  -4.87%  "

# Row 47: update D47, E47
$ws.Range("D47").Formula = "'3.52"
$ws.Range("E47").Value = "  -0.68%  "

# Row 48: update D48, E48
$ws.Range("D48").Formula = "'20.07"
$ws.Range("E48").Value = "  -3.01%  "

# Row 49: update E49
$ws.Range("E49").Value = "  -2.55%  "

# Row 50: update D50, E50
$ws.Range("D50").Formula = "'0.571"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51: update D51, E51
$ws.Range("D51").Formula = "'0.0921"
$ws.Range("E51").Value = "  -1.48%  "
